# "Generate Report for Handback"
# Updates the localization-status report with handback results:
#  - Status flips from "Ready for handoff" to "Handback transform failed"
#    on the Overview sheet (zh-cn/de-de columns) and on the per-locale
#    "Status" column of the zh-cn / de-de sheets.
#  - The per-locale "Error Detail" column (P) gets the failure detail
#    message, and is widened to fit it.

$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

$newStatus = "Handback transform failed"

# Overview sheet: zh-cn (E3) and de-de (F3) status cells for the
# b5512d5a-... row.
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# zh-cn sheet: Status column (C) for the same row, plus the Error Detail
# (P) message explaining the handback transform failure.
$wsZhCn.Range("C3").Value = $newStatus
$wsZhCn.Range("P3").Value = "Handback file name: nqlh2vdi.jaf is different with handoff file name: b5512d5a-de2d-4c1c-a4df-ff538edf848d.906bae67b670507c91213cdaf7f95eead7b3417f.zh-cn."
$wsZhCn.Columns.Item(16).ColumnWidth = 39.166666666666664

# de-de sheet: same pair of edits.
$wsDeDe.Range("C3").Value = $newStatus
$wsDeDe.Range("P3").Value = "Handback file name: nqlh2vdi.jaf is different with handoff file name: b5512d5a-de2d-4c1c-a4df-ff538edf848d.906bae67b670507c91213cdaf7f95eead7b3417f.de-de."
$wsDeDe.Columns.Item(16).ColumnWidth = 39.166666666666664
